# [Kadastro App] Yeni kayit eklendi: 3010
# Appends the new record (Kayit No 3010) as the next row on both the
# master "Kayitlar" sheet and the district-filtered "Erdemli" sheet,
# mirroring the existing rows exactly (all values stored as text).

$wb = $excel.ActiveWorkbook

$newRow = @("3010", "2025-09-11", "Erdemli", "1", "CİNS DEĞ.", "AYHAN KARADAYI (K.Teknisyeni), EMİNE ALANLI KIRCILI (K.Mühendisi)")

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Existing data occupies rows 1..68 (header + 67 records); the new
    # record lands on row 69, directly under the last existing row.
    $targetRow = 69

    $rng = $ws.Range("A$targetRow`:F$targetRow")
    # Force text formatting first so numeric-looking values ("3010", "1")
    # and the date-looking value ("2025-09-11") are stored as literal text
    # strings, same as every other cell in these columns.
    $rng.NumberFormat = "@"

    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $newRow[$col - 1]
    }

    # Drop back to the default "Normal" style so the new row doesn't carry
    # an extra explicit number-format style compared to the surrounding rows.
    $rng.Style = "Normal"
}
